# Correspondance.xlsx edit
# "Fin bâtiments + tous twy/apron / reste arbres + divers objets"
#
# Content changes (sheet "Feuil1"):
#  - Row 28 (S8): rename correspondence text for "Aero club Ouest Parisien"
#      -> "Aero club Ouest Parisien - ACOP"
#  - Row 29 (S9): add new correspondence text "JC Decaux"
#  - Row 30 (S10): split/rename "HélixAero / Clean Aéro Service"
#      -> "Clean Aéro Service"
#  - Row 31 (S11): numeric code 243 -> 245, add new correspondence text
#      "Restaurant Air & Cook"
#  - Row 40 (W4): rename "IX AIR / First Flight"
#      -> "IX AIR Helix Aero  First Flight"
#
# Plus cosmetic window/selection state: window minimized, active selection
# moved from C32 to C11 (no more frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Order chosen to reproduce the shared-string table append order of the
# original edit (new strings are appended in this sequence).
$ws.Range("C31").Value = "Restaurant Air & Cook"
$ws.Range("C28").Value = "Aero club Ouest Parisien - ACOP"
$ws.Range("C30").Value = "Clean Aéro Service"
$ws.Range("C40").Value = "IX AIR Helix Aero  First Flight"
$ws.Range("C29").Value = "JC Decaux"

$ws.Range("B31").Value = 245

# View / selection state
$ws.Range("C11").Select()
$wb.Windows.Item(1).WindowState = -4140
